$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '59.812.09'
$ws.Range("E2").Value = '  +1.21%  '
$ws.Range("D3").Value = '2.650.20'
$ws.Range("E3").Value = '  +2.34%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '536.64'
$ws.Range("E5").Value = '  +1.18%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '145.41'
$ws.Range("E6").Value = '  +3.97%  '
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("E8").Value = '  +1.09%  '
$ws.Range("D9").Value = '2.667.86'
$ws.Range("E9").Value = '  +2.55%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.73'
$ws.Range("E10").Value = '  +4.52%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.103'
$ws.Range("E11").Value = '  +1.34%  '
$ws.Range("E12").Value = '  +1.43%  '
$ws.Range("E13").Value = '  -1.11%  '
$ws.Range("D14").Value = '3.124.47'
$ws.Range("E14").Value = '  +2.39%  '
$ws.Range("D15").Value = '59.743.86'
$ws.Range("E15").Value = '  +1.16%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '21.27'
$ws.Range("E16").Value = '  +4.30%  '
$ws.Range("D17").Value = '2.647.98'
$ws.Range("E17").Value = '  +2.38%  '
$ws.Range("E18").Value = '  +1.39%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '344.91'
$ws.Range("E19").Value = '  -0.35%  '
$ws.Range("E20").Value = '  +2.27%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.24'
$ws.Range("E21").Value = '  +1.48%  '
$ws.Range("E22").Value = '  -0.51%  '
$ws.Range("E23").Value = '  +0.03%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '66.83'
$ws.Range("E24").Value = '  -0.86%  '
$ws.Range("E25").Value = '  +2.44%  '
$ws.Range("E27").Value = '  +0.07%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.31'
$ws.Range("E28").Value = '  +2.28%  '
$ws.Range("D29").Value = '0.0₃0750'
$ws.Range("E29").Value = '  +2.62%  '
$ws.Range("E30").Value = '  -0.07%  '
$ws.Range("E31").Value = '  +2.81%  '
$ws.Range("B32").Value = 'Aptos'
$ws.Range("C32").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.87'
$ws.Range("E32").Value = '  +1.02%  '
$ws.Range("B33").Value = 'EthereumClassic'
$ws.Range("C33").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '19.06'
$ws.Range("E33").Value = '  +1.69%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '150.19'
$ws.Range("E34").Value = '  +0.50%  '
$ws.Range("E35").Value = '  +1.79%  '
$ws.Range("E36").Value = '  +3.27%  '
$ws.Range("E37").Value = '  +0.13%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.841'
$ws.Range("E38").Value = '  +1.72%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '296.11'
$ws.Range("E39").Value = '  +9.15%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.823'
$ws.Range("E40").Value = '  -0.42%  '
$ws.Range("E41").Value = '  +2.13%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.999'
$ws.Range("E42").Value = '  +0.09%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.605'
$ws.Range("E43").Value = '  +1.83%  '
$ws.Range("E44").Value = '  +5.35%  '
$ws.Range("B45").Value = 'WhiteBITCoin'
$ws.Range("C45").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '10.74'
$ws.Range("E45").Value = '  -0.08%  '
$ws.Range("B46").Value = 'EnergySwap'
$ws.Range("C46").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '19.34'
$ws.Range("E46").Value = '  +5.15%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0954'
$ws.Range("E47").Value = '  -0.45%  '
$ws.Range("E48").Value = '  +2.45%  '
$ws.Range("D49").Value = '1.970.74'
$ws.Range("E49").Value = '  +1.27%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '4.57'
$ws.Range("E50").Value = '  +2.65%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '18.38'
$ws.Range("E51").Value = '  +1.26%  '
